# Weekly fruit/vegetable price refresh: the "Fecha" (date) + price/volume
# columns for each data row get re-drawn from another row in the same
# subset (a full reshuffle of rows 2-33 across columns D, J, K, L, M, N,
# O, P, Q), while the descriptive columns (A, B, C, E, F, G, H, I, R) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 33

# Columns that move together as a unit per source row.
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot every "movable" cell before writing anything back, so the
# permutation below reads only original values (never a value already
# overwritten earlier in the loop).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2()
    }
    $snapshot[$r] = $rowVals
}

# Target row -> source row (the row whose original D/J/K/L/M/N/O/P/Q
# values now populate the target row).
$perm = @{
    2 = 8;   3 = 29;  4 = 6;   5 = 22;  6 = 28;  7 = 25;  8 = 10;  9 = 24;
    10 = 14; 11 = 17; 12 = 30; 13 = 23; 14 = 31; 15 = 19; 16 = 15; 17 = 4;
    18 = 13; 19 = 3;  20 = 32; 21 = 9;  22 = 16; 23 = 21; 24 = 33; 25 = 26;
    26 = 12; 27 = 7;  28 = 5;  29 = 20; 30 = 11; 31 = 2;  32 = 18; 33 = 27
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $perm[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
